# The presentation ships with two themes:
#   theme1.xml -> bound to the (only) Slide Master, originally the
#                 "Integral" / "Red Violet" design.
#   theme2.xml -> bound to the Notes Master, originally the stock
#                 "Office Theme" / "Office" color scheme.
#
# The target revision swaps which design is "active": the Slide Master
# (theme1.xml) is repainted with the plain built-in "Office" color
# scheme (dk1=000000, lt1=FFFFFF, dk2=44546A, lt2=E7E6E6,
# accent1-6=5B9BD5/ED7D31/A5A5A5/FFC000/4472C4/70AD47,
# hlink=0563C1, folHlink=954F72) - i.e. exactly what used to live in
# theme2.xml. We reproduce that through the live PowerPoint color-scheme
# object model: SlideMaster.ColorScheme.Colors(n) maps 1:1 onto the 12
# OOXML clrScheme slots in document order (dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink), and writing RGBColor.RGB rewrites
# the corresponding <a:srgbClr val="…"/> in the slide master's theme
# part.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$cs = $master.ColorScheme

# RGB() helper: PowerPoint COM colors are packed 0xBBGGRR (little endian
# red-green-blue), NOT 0xRRGGBB - build the values from hex bytes so the
# intent stays readable.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# New ("Office") scheme, applied in clrScheme slot order 1-12.
$cs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$cs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$cs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$cs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$cs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$cs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$cs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$cs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$cs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$cs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$cs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$cs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
